$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) Insert a new column at Z (column 26), shifting web/webalert/webcookie/ws/ws.async/xml
#    one column to the right (Z->AA, AA->AB, AB->AC, AC->AD, AD->AE, AE->AF).
$ws.Columns.Item(26).Insert()

# 2) Populate the newly inserted Z column with the new "tn.5250" category
#    (header in row 1, 5 data rows below it).
$ws.Range("Z1").Value = "tn.5250"
$ws.Range("Z2").Value = "close(profile)"
$ws.Range("Z3").Value = "open(profile)"
$ws.Range("Z4").Value = "saveText(profile,var)"
$ws.Range("Z5").Value = "typeKeys(profile,keystrokes)"
$ws.Range("Z6").Value = "updateScreenFields(profile)"

# 3) Update the "image" category (column K): rename colorbit's first parameter,
#    add the new ocr(image,saveVar) command, and re-sort alphabetically.
$ws.Range("K2").Value = "colorbit(image,bit,saveTo)"
$ws.Range("K6").Value = "ocr(image,saveVar)"
$ws.Range("K7").Value = "resize(image,width,height,saveTo)"
$ws.Range("K8").Value = "saveDiff(var,baseline,actual)"

# 4) Update the "target" category (column A) to insert "tn.5250" alphabetically
#    between "step" and "web".
$ws.Range("A26").Value = "tn.5250"
$ws.Range("A27").Value = "web"
$ws.Range("A28").Value = "webalert"
$ws.Range("A29").Value = "webcookie"
$ws.Range("A30").Value = "ws"
$ws.Range("A31").Value = "ws.async"
$ws.Range("A32").Value = "xml"

# 5) Touch the new rightmost column (AG) so the sheet's recorded dimension grows
#    by one column, matching the pre-existing +1 padding convention of this sheet.
$ws.Range("AG144").NumberFormat = "General"

# 6) Fix up the defined names so they point at the correct (shifted) ranges
#    and reflect the grown "image"/"target" lists, plus register "tn.5250".
$wb.Names.Item("image").RefersTo = "='#system'!`$K`$2:`$K`$8"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$144"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"
$wb.Names.Add("tn.5250", "='#system'!`$Z`$2:`$Z`$6")
